$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) are stored as plain text in this sheet,
# so force text format before writing to avoid Excel auto-converting
# numeric-looking strings (e.g. '1.001') into numbers.

# --- Update Price (D) / Volume(1h) (E) columns for rows whose coin stayed the same ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.871.31"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.757.15"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.60"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.98%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4578"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3498"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.46%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "41.93"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07358"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.085"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.59"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.980"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.171"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.758.29"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.67"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001053"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06425"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.83"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.755"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.888.77"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.09%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.44%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +4.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.86"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.05"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.960.39"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.155"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.48"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.073"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09244"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.659"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.546"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.72"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.10%  "

# --- Rows 36 & 37: Hedera/VeChain swapped rank order ---
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02267"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.78%  "

$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06096"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.28%  "

# --- Rows 38-48: continue updating Price / Volume(1h) ---
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2060"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.900"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6177"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.177"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.365"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.766"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.15"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.726"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5789"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "123.34"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.927"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.61%  "

# --- Rows 49 & 50: EOS/Cronos swapped rank order ---
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06796"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.42%  "

$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.122"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.77%  "

# --- Row 51 ---
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.20"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.13%  "

